# Update view-count figures (column F) for several rows across two sheets
# as reflected in the "output generated at 456a3b4" data refresh.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 258
$wsExhibit.Range("F8").Value = 1941
$wsExhibit.Range("F10").Value = 4550

# Sheet "全部类型" (All types) - same events duplicated with different row numbers
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 258
$wsAll.Range("F12").Value = 1941
$wsAll.Range("F14").Value = 4550
